# 5.4.1.xlsx — add a 2020 data column (E) alongside the existing 2015 column (D).
#
# Steps:
#   1. Copy the formatting of column D (rows 3-8) into column E so the new
#      column inherits the same borders / fonts / alignment as the 2015 one.
#   2. Write the 2020 figures into E4:E8.
#   3. Give E7 a one-decimal "0.0" custom number format (this is what creates
#      the new numFmt/cellXfs entries in styles.xml).
#   4. Leave the active selection on B15, matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clone column D's formatting (borders/fonts/alignment) onto column E for
#    the data block (rows 3-8); no values are copied by this paste.
$ws.Range("D3:D8").Copy()
$ws.Range("E3:E8").PasteSpecial(-4122)  # xlPasteFormats

# 2. Fill in the 2020 values.
$ws.Range("E4").Value = 2020
$ws.Range("E5").Value = 11.5
$ws.Range("E6").Value = 2.6
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 0.3

# 3. E7 gets its own one-decimal display format (adds numFmt 164 = "0.0").
$ws.Range("E7").NumberFormat = "0.0"

# 4. Restore the saved selection/active cell.
$ws.Range("B15").Select()
